# Update layout master barang multi satuan dan tambah master barang pelanggan
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update values: N2:N5 3000 -> 19000
$ws.Range("N2").Value = 19000
$ws.Range("N3").Value = 19000
$ws.Range("N4").Value = 19000
$ws.Range("N5").Value = 19000

# Update value G5: 200 -> 2000
$ws.Range("G5").Value = 2000

# Update the active selection on the sheet to N8
$ws.Activate()
$ws.Range("N8").Select()
